# Re-order the elements printed for each First(...) set in column D so that
# the textual representation matches the new (nested-array-initializer-aware)
# set iteration order produced by the CFG first-set generator.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '{''eklabool'', ''andamhie'', ''shimenet'', ''chika'', ''anda'', ''naur''}'
$ws.Range("D3").Value = '{''eklabool'', ''andamhie'', ''shimenet'', ''chika'', ''anda'', ''naur'', ''λ''}'
$ws.Range("D4").Value = '{''eklabool'', ''andamhie'', ''chika'', ''shimenet'', ''anda'', ''naur'', ''λ''}'
$ws.Range("D6").Value = '{''eklabool'', ''anda'', ''andamhie'', ''chika'', ''shimenet''}'
$ws.Range("D7").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D8").Value = '{''('', ''λ'', ''['', ''=''}'
$ws.Range("D10").Value = '{''eklabool'', ''andamhie'', ''anda'', ''chika'', ''λ''}'
$ws.Range("D11").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D13").Value = '{''λ'', ''[''}'
$ws.Range("D14").Value = '{''λ'', ''[''}'
$ws.Range("D15").Value = '{''λ'', ''[''}'
$ws.Range("D16").Value = '{''λ'', ''[''}'
$ws.Range("D17").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D18").Value = '{''andamhie_literal'', ''chika_literal'', ''korik'', ''anda_literal'', ''eme''}'
$ws.Range("D19").Value = '{''eme'', ''korik''}'
$ws.Range("D20").Value = '{''('', ''λ'', ''[''}'
$ws.Range("D21").Value = '{''('', ''λ'', ''[''}'
$ws.Range("D26").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D28").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''{'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D30").Value = '{''eklabool'', ''anda'', ''andamhie'', ''shimenet'', ''chika'', ''λ''}'
$ws.Range("D31").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''serve'', ''forda'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D32").Value = '{''eklabool'', ''anda'', ''naur'', ''andamhie'', ''chika'', ''λ''}'
$ws.Range("D33").Value = '{''eklabool'', ''andamhie'', ''chika'', ''anda'', ''naur''}'
$ws.Range("D34").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''λ'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D36").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D37").Value = '{''λ'', ''id''}'
$ws.Range("D38").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''serve'', ''forda'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D39").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D40").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D42").Value = '{''%'', ''!='', ''**'', ''>'', ''<='', ''&&'', ''=='', ''<'', ''/'', ''>='', ''||'', ''-'', ''λ'', ''*'', ''+'', ''//''}'
$ws.Range("D43").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D44").Value = '{''!'', ''λ'', ''-''}'
$ws.Range("D45").Value = '{''andamhie_literal'', ''len'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''('', ''eme''}'
$ws.Range("D47").Value = '{''--'', ''('', ''λ'', ''['', ''++''}'
$ws.Range("D48").Value = '{''--'', ''++''}'
$ws.Range("D49").Value = '{''andamhie_literal'', ''chika_literal'', ''korik'', ''anda_literal'', ''eme''}'
$ws.Range("D50").Value = '{''%'', ''!='', ''**'', ''>'', ''<='', ''&&'', ''=='', ''<'', ''/'', ''>='', ''||'', ''-'', ''*'', ''+'', ''//''}'
$ws.Range("D51").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''serve'', ''forda'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D52").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''gogogo'', ''serve'', ''forda'', ''amaccana'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D53").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''gogogo'', ''serve'', ''forda'', ''amaccana'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D54").Value = '{''**='', ''//='', ''%='', ''['', ''/='', ''-='', ''+='', ''='', ''('', ''*=''}'
$ws.Range("D55").Value = '{''+='', ''='', ''**='', ''//='', ''%='', ''/='', ''*='', ''-=''}'
$ws.Range("D56").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''{'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D57").Value = '{''eklabool'', ''andamhie'', ''anda'', ''id'', ''chika''}'
$ws.Range("D58").Value = '{''eklabool'', ''andamhie'', ''anda'', ''chika'', ''λ''}'
$ws.Range("D59").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D63").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D65").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D66").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''serve'', ''forda'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D68").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''gogogo'', ''serve'', ''forda'', ''amaccana'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D69").Value = '{''λ'', ''ganern''}'
$ws.Range("D70").Value = '{''λ'', ''ganern''}'
$ws.Range("D73").Value = '{''eklabool'', ''andamhie'', ''anda'', ''chika'', ''λ''}'
$ws.Range("D74").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D75").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D76").Value = '{''λ'', ''step''}'
$ws.Range("D77").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''λ'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D78").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D79").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''gogogo'', ''serve'', ''forda'', ''amaccana'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D82").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''gogogo'', ''serve'', ''forda'', ''amaccana'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D84").Value = '{''betsung'', ''λ''}'
$ws.Range("D85").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D86").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''serve'', ''forda'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D88").Value = '{''betsung'', ''λ''}'
$ws.Range("D89").Value = '{''eklabool'', ''id'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''adele'', ''andamhie'', ''--'', ''gogogo'', ''serve'', ''forda'', ''amaccana'', ''adelete'', ''versa'', ''chika'', ''λ'', ''pak''}'
$ws.Range("D91").Value = '{''λ'', ''ditech''}'
$ws.Range("D92").Value = '{''amaccana'', ''λ'', ''gogogo''}'
$ws.Range("D93").Value = '{''λ'', ''push''}'
$ws.Range("D94").Value = '{''--'', ''++'', ''id''}'
$ws.Range("D95").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''λ'', ''{'', ''korik'', ''-'', ''('', ''eme''}'
